$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new data row (2007) is now the first data row, shifting every existing
# row down by one and adding a brand new last row (2024 -> 2025 forecast).
# Apply the date-column style (used by A2:A18 already) to the new last row,
# A19, so formatting stays consistent with the rest of the column.
$ws.Range("A18").Copy()
$ws.Range("A19").PasteSpecial(-4122) # xlPasteFormats

$data = @(
    @(39400, 2007, 7.226520411029047, 2008, 8.987952903094421),
    @(39765, 2008, 4.268860212333636, 2009, 6.77211531652997),
    @(40130, 2009, -7.266312015249776, 2010, -1.420242831007679),
    @(40494, 2010, 6.958243460951929, 2011, 8.31992099499319),
    @(40862, 2011, 9.469137444079934, 2012, 8.571528775834981),
    @(41228, 2012, 3.358206407534947, 2013, 5.745831525574441),
    @(41592, 2013, 0.3081076735359067, 2014, 2.847379875994704),
    @(41957, 2014, 3.901355411819707, 2015, 3.690459963535009),
    @(42321, 2015, 5.331683351557981, 2016, 4.636575318346536),
    @(42689, 2016, 3.254758369308375, 2017, 2.76788332063731),
    @(43053, 2017, 5.246209615995667, 2018, 4.659862065074982),
    @(43418, 2018, 4.86255966374296, 2019, 4.636196713604357),
    @(43783, 2019, 2.764740011159428, 2020, 2.471557257221946),
    @(44159, 2020, -7.260793671746435, 2021, -1.387795042833839),
    @(44525, 2021, 4.097586525396268, 2022, 2.450242954096926),
    @(44890, 2022, 7.824284864703746, 2023, 4.834990656989402),
    @(45254, 2023, -1.24502235313334, 2024, -2.488220481262082),
    @(45618, 2024, -1.735114423676209, 2025, 1.056286187957367)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}
